# Re-case a handful of header cells on row 1 of Sheet1 to title-case
# (every word capitalized) instead of only-first-word-capitalized, and
# move the active selection from AR8 (with topLeftCell AP1) to O6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerUpdates = @{
    "A1"  = "Local Authority"
    "C1"  = "Provider Type"
    "D1"  = "Provider Subtype"
    "F1"  = "Setting Name"
    "G1"  = "Registration Date"
    "H1"  = "Registration Status"
    "I1"  = "Deregistration Date"
    "J1"  = "Placement Code"
    "O1"  = "Setting Address Postcode"
    "AA1" = "Placement Provider Code"
}

foreach ($cellRef in $headerUpdates.Keys) {
    $ws.Range($cellRef).Value = $headerUpdates[$cellRef]
}

# Update the visible selection to match the saved view (no frozen/scrolled
# topLeftCell override, single-cell selection at O6).
$ws.Range("O6").Select()
